$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 1 with two runs of 1,2,3 across A1:F1
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 3
$ws.Range("D1").Value = 1
$ws.Range("E1").Value = 2
$ws.Range("F1").Value = 3

# Leave the final selection on F1, matching the last-edited cell
$ws.Range("F1").Select()
